{"js": "// Office.js (Word JS API) edit script.\n// Applies the Review_480 edit: updates the date, swaps the paper title and\n// review body paragraphs from the \"Procedural Knowledge\" review to the\n// \"Tax Evasion\" review, and appends a closing remark + new arXiv link\n// (preceded by a line break), replacing the old link paragraph.\n\nconst newTexts = [\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7: 04.07.25\",\n  \"Investigating Tax Evasion Emergence Using Dual Large Language Model and Deep Reinforcement Learning Powered Agent-based Simulation\",\n  \"\u05ea\u05e4\u05e0\u05d9\u05ea \u05de\u05e4\u05ea\u05d9\u05e2\u05d4 \u05de\u05ea\u05e8\u05d7\u05e9\u05ea \u05d1\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1-LLMs \u05d1\u05ea\u05d7\u05d5\u05de\u05d9\u05dd \\\"\u05e8\u05db\u05d9\u05dd\\\" \u05d9\u05d5\u05ea\u05e8 \u05db\u05de\u05d5 \u05e4\u05e1\u05d9\u05db\u05d5\u05dc\u05d5\u05d2\u05d9\u05d4, \u05e1\u05d5\u05e6\u05d9\u05d5\u05dc\u05d5\u05d2\u05d9\u05d4 \u05d5\u05d0\u05e4\u05d9\u05dc\u05d5 \u05db\u05dc\u05db\u05dc\u05d4. \u05d0\u05d5\u05de\u05e0\u05dd LLMs \u05dc\u05d0 \\\"\u05d7\u05d5\u05e9\u05d1\u05d9\u05dd\\\" \u05db\u05de\u05d5 \u05d1\u05e0\u05d9 \u05d0\u05d3\u05dd \u05d1\u05e8\u05de\u05ea \u05d4\u05d0\u05d9\u05e0\u05d3\u05d9\u05d1\u05d9\u05d3\u05d5\u05d0\u05dc, \u05d0\u05d1\u05dc \u05de\u05e1\u05ea\u05d1\u05e8 \u05e9\u05d4\u05dd \u05db\u05d1\u05e8 \u05de\u05d7\u05db\u05d9\u05dd \u05d0\u05ea \u05d0\u05d9\u05da \u05e9\u05d0\u05e0\u05d7\u05e0\u05d5 \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05db\u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d4. \",\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05e0\u05d5\u05db\u05d7\u05d9 \u05de\u05e6\u05d9\u05d2 \u05d2\u05d9\u05e9\u05d4 \u05d7\u05d3\u05e9\u05e0\u05d9\u05ea \u05dc\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1 LLMs \u05dc\u05d7\u05e7\u05e8 \u05d4\u05ea\u05d7\u05de\u05e7\u05d5\u05ea \u05de\u05de\u05e1 \u05d1\u05e2\u05d6\u05e8\u05ea \u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4. \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05d4\u05e0\u05d9\u05d7 \u05de\u05e8\u05d0\u05e9 \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05e9\u05dc \u05d4\u05e2\u05dc\u05de\u05ea \u05de\u05e1, \u05db\u05de\u05d5 \u05e9\u05e2\u05e9\u05d5 \u05db\u05dc \u05d4\u05de\u05d7\u05e7\u05e8\u05d9\u05dd \u05e2\u05d3 \u05d4\u05d9\u05d5\u05dd, \u05d4\u05de\u05d7\u05e7\u05e8 \u05de\u05ea\u05de\u05e7\u05d3 \u05d1\u05d4\u05d5\u05e4\u05e2\u05d4 \u05d5\u05d1\u05d3\u05d9\u05e0\u05de\u05d9\u05e7\u05d4 \u05e9\u05dc \u05ea\u05d5\u05e4\u05e2\u05d4 \u05d6\u05d5 \u05d1\u05e7\u05e8\u05d1 \u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d9\u05d4. \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e1\u05d5\u05db\u05e0\u05d9\u05dd (agent-based simulation) \u05d4\u05de\u05e9\u05dc\u05d1\u05ea LLMs \u05d5\u05dc\u05de\u05d9\u05d3\u05d4 \u05e2\u05de\u05d5\u05e7\u05d4 \u05e2\u05dd \u05d7\u05d9\u05d6\u05d5\u05e7\u05d9\u05dd (deep reinforcement learning), \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05d1\u05d5\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05d4\u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05d9\u05d5\u05ea \u05db\u05dc\u05db\u05dc\u05d9\u05d5\u05ea \u05d1\u05dc\u05ea\u05d9 \u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05d5\u05ea (\u05de\u05d4 \u05e9\u05d4\u05e8\u05d1 \u05de\u05db\u05d9\u05e8\u05d9\u05dd \u05d1\u05ea\u05d5\u05e8 \\\"\u05db\u05dc\u05db\u05dc\u05d4 \u05e9\u05d7\u05d5\u05e8\u05d4\\\") \u05dc\u05d4\u05d2\u05d9\u05d7 \u05d1\u05d0\u05d5\u05e4\u05df \u05e1\u05e4\u05d5\u05e0\u05d8\u05e0\u05d9, \u05d5\u05dc\u05d0 \u05db\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05de\u05d5\u05d2\u05d3\u05e8\u05ea \u05de\u05e8\u05d0\u05e9. \u05e2\u05dc \u05d0\u05e3 \u05d4\u05ea\u05e8\u05d5\u05de\u05d4 \u05e9\u05dc\u05d5 \u05dc\u05db\u05dc\u05db\u05dc\u05d4, \u05de\u05d4 \u05e9\u05dc\u05db\u05e0\u05e8\u05d0\u05d4 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d0\u05d5\u05ea\u05e0\u05d5 \u05d6\u05d4 \u05d4\u05e9\u05d9\u05de\u05d5\u05e9 \u05d4\u05d9\u05d9\u05d7\u05d5\u05d3\u05d9 \u05d1 LLMs \u05d5\u05d1 DRL \u05db\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc \u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05e9\u05d9\u05d5\u05d3\u05e2 \u05dc\u05e7\u05d7\u05ea \u05de\u05d9\u05d3\u05e2 \u05dc\u05d0 \u05e4\u05d5\u05e8\u05de\u05dc\u05d9 (\u05db\u05de\u05d5 \u05ea\u05d9\u05d0\u05d5\u05e8 \u05d0\u05d9\u05e9\u05d9\u05d5\u05ea) \u05d5\u05dc\u05d4\u05db\u05e0\u05d9\u05e1 \u05d0\u05d5\u05ea\u05d5 \u05dc\u05e1\u05d8 \u05d4\u05e9\u05d9\u05e7\u05d5\u05dc\u05d9\u05dd \u05d4\u05e4\u05d5\u05e8\u05de\u05dc\u05d9 \u05e9\u05dc \u05e1\u05d5\u05db\u05df - \u05dc\u05de\u05e9\u05dc \u05db\u05de\u05d4 \u05d0\u05d5\u05e4\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc \u05e4\u05e8\u05e1\u05d5\u05e0\u05dc\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d1\u05de\u05e2\u05e8\u05db\u05d5\u05ea \u05d0\u05e4\u05e9\u05e8 \u05dc\u05e2\u05e9\u05d5\u05ea \u05e2\u05dc \u05d2\u05d1\u05d9 \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05d6\u05d4.\",\n  \"\u05d4\u05e9\u05d9\u05d8\u05d4 \u05d1\u05de\u05d0\u05de\u05e8 \u05de\u05ea\u05e2\u05de\u05e7\u05ea \u05d1\u05d9\u05e6\u05d9\u05e8\u05ea \u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05d4\u05de\u05d3\u05de\u05d4 \u05db\u05dc\u05db\u05dc\u05d4 \u05e1\u05d2\u05d5\u05e8\u05d4 (\u05d9\u05e9 \u05de\u05e1\u05d7\u05e8 \u05d1\u05d9\u05df \u05d0\u05e0\u05e9\u05d9\u05dd \u05d1\u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d4 \u05e9\u05de\u05e0\u05e1\u05d4 \u05dc\u05d4\u05d8\u05d9\u05d9\u05d1 \u05e2\u05dd \u05de\u05e6\u05d1\u05dd), \u05e9\u05d1\u05d4 \u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05e4\u05d5\u05e2\u05dc\u05d9\u05dd \u05d5\u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea. \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05d4\u05e9\u05ea\u05de\u05e9\u05d5 \u05d1\u05de\u05d1\u05e0\u05d4 \u05db\u05dc\u05db\u05dc\u05d4 \u05d3\u05d5\u05de\u05d4 \u05dc\u05d6\u05d0\u05ea \u05e9\u05dc \u05d0\u05e8\u05e6\u05d5\u05ea \u05d4\u05d1\u05e8\u05d9\u05ea \u05e2\u05dd \u05d3\u05d9\u05d5\u05d5\u05d7 \u05de\u05e1 \u05e2\u05e6\u05de\u05d9 \u05db\u05d3\u05d9 \u05dc\u05d0\u05e4\u05e9\u05e8 \u05dc\u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05d4\u05d6\u05d3\u05de\u05e0\u05d5\u05ea \u05dc\u05d4\u05e2\u05dc\u05d9\u05dd \u05de\u05e1 \u05db\u05d4\u05d7\u05dc\u05d8\u05d4. \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05d9\u05e9 \u05d4\u05de\u05d5\u05df \u05e1\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05de\u05e1 \u05d5\u05d6\u05d4 \u05de\u05e1\u05ea\u05d1\u05da \u05de\u05d4\u05e8 (\u05ea\u05e9\u05d0\u05dc\u05d5 \u05d0\u05ea \u05e8\u05d5\u05d0\u05d4 \u05d4\u05d7\u05e9\u05d1\u05d5\u05df \u05e9\u05dc\u05db\u05dd) \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05d4\u05ea\u05e8\u05db\u05d6\u05d5 \u05d1\u05e9\u05e0\u05d9 \u05e1\u05d5\u05d2\u05d9 \u05de\u05e1 - \u05de\u05e1 \u05d4\u05db\u05e0\u05e1\u05d4 \u05d5\u05de\u05e1 \u05e2\u05e8\u05da \u05de\u05d5\u05e1\u05e3 (\u05de\u05e2\\\"\u05de). \",\n  \"\u05d4\u05dd \u05d2\u05dd \u05e1\u05d9\u05de\u05dc\u05e6\u05d5 \u05e8\u05e9\u05d5\u05d9\u05d5\u05ea \u05d4\u05d7\u05d5\u05e7 \u05d5\u05ea\u05d5\u05e2\u05dc\u05ea \u05e9\u05d4\u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05de\u05d4\u05de\u05d3\u05d9\u05e0\u05d4 \u05de\u05d4\u05de\u05e1 \u05e9\u05d4\u05dd \u05de\u05e9\u05dc\u05de\u05d9\u05dd \u05dc\u05d4. \u05d4\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05e2\u05e6\u05de\u05d4 \u05de\u05d4\u05d5\u05d5\u05d4 \u05ea\u05e9\u05ea\u05d9\u05ea \u05dc\u05d7\u05dc\u05e7 \u05d4\u05de\u05e8\u05db\u05d6\u05d9 \u05e9\u05dc \u05d4\u05e2\u05d1\u05d5\u05d3\u05d4 - \u05de\u05d5\u05d3\u05dc \u05e7\u05d1\u05dc\u05ea \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05e9\u05dc \u05d4\u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05d1\u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d4. \u05db\u05d3\u05d9 \u05dc\u05d9\u05d9\u05e6\u05e8 \u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d4 \u05d4\u05d8\u05e8\u05d5\u05d2\u05e0\u05d9\u05ea, \u05de\u05d5\u05d7 \u05e9\u05dc \u05e1\u05d5\u05db\u05df \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc LLM \u05d5 DRL \u05db\u05d0\u05e9\u05e8 LLMs \u05e9\u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05ea\u05d9\u05d0\u05d5\u05e8 \u05e9\u05dc \u05d4\u05d0\u05d5\u05e4\u05d9 \u05e9\u05dc \u05d4\u05e1\u05d5\u05db\u05df (\u05de\u05d1\u05d5\u05e1\u05e1 \u05e2\u05dc \u05d4\u05ea\u05d5\u05db\u05df \u05e9\u05d4\u05d5\u05d0 \u05de\u05e4\u05e8\u05e1\u05dd \u05d1\u05d8\u05d5\u05d5\u05d9\u05d8\u05e8 \u05dc\u05de\u05e9\u05dc), \u05d4\u05d9\u05e1\u05d8\u05d5\u05e8\u05d9\u05d9\u05ea \u05d4\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05db\u05d8\u05e7\u05e1\u05d8, \u05d5\u05d0\u05ea \u05db\u05dc \u05d4\u05de\u05d9\u05d3\u05e2 \u05e2\u05dc \u05d4\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05db\u05d8\u05e7\u05e1\u05d8 \u05d2\u05dd \u05db\u05df. \",\n  \"\u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4-context prompt \u05d4\u05d6\u05d4 \u05d4 LLM \u05e0\u05e9\u05d0\u05dc \\\"\u05db\u05de\u05d4 \u05de\u05e1 \u05d0\u05e0\u05d9 \u05e6\u05e8\u05d9\u05da \u05dc\u05e9\u05dc\u05dd?\\\". \u05d4\u05de\u05e1\u05e4\u05e8 \u05e9\u05d4 LLM \u05de\u05d7\u05d6\u05d9\u05e8, \u05de\u05d5\u05e2\u05d1\u05e8 \u05db\u05e7\u05dc\u05d8 \u05dc\u05de\u05d5\u05d3\u05dc  DRL \u05e9\u05de\u05e7\u05d1\u05dc \u05d2\u05dd \u05d0\u05ea \u05e9\u05d0\u05e8 \u05d4\u05d3\u05d0\u05d8\u05d4 \u05e9\u05e7\u05d9\u05d1\u05dc \u05d4 LLM \u05d0\u05d1\u05dc \u05d2\u05dd \u05db\u05de\u05d4 \u05d4\u05e1\u05d5\u05db\u05df \\\"\u05d4\u05e8\u05e4\u05ea\u05e7\u05e0\u05d9\\\" \u05db\u05e4\u05e8\u05de\u05d8\u05e8 \u05e9\u05d4-DRL \u05de\u05e9\u05ea\u05de\u05e9 \u05db\u05d3\u05d9 \u05dc\u05e2\u05e9\u05d5\u05ea \u05d0\u05e7\u05e1\u05e4\u05dc\u05d5\u05e8\u05e6\u05d9\u05d4. \u05db\u05dc\u05d5\u05de\u05e8, \u05d4 LLM \u05de\u05d7\u05d6\u05d9\u05e8 \u05d4\u05d7\u05dc\u05d8\u05d4 \u05e8\u05d0\u05e9\u05d5\u05e0\u05d9\u05ea \u05e9\u05d0\u05d5\u05ea\u05d4, \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05dc\u05e7\u05d8 \u05e9\u05dc \u05d4 LLM \u05d5\u05e2\u05d5\u05d3 \u05de\u05e9\u05ea\u05e0\u05d4 exploration (\u05d0\u05d4\u05d1\u05ea \u05e1\u05d9\u05db\u05d5\u05df \u05d1\u05de\u05d5\u05d1\u05df \u05d4\u05db\u05dc\u05db\u05dc\u05d9) \u05de\u05e7\u05d1\u05dc \u05d2\u05dd DRL \u05e9\u05de\u05e7\u05d1\u05dc \u05d4\u05d7\u05dc\u05d8\u05d4 \u05d1\u05e2\u05e6\u05de\u05d5 \u05e9\u05d4\u05d9\u05e0\u05d4 \u05d2\u05dd \u05d4\u05e1\u05d5\u05e4\u05d9\u05ea. \",\n  \"\u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05dc\u05de\u05d9\u05d3\u05d4 \u05d4\u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05e2\u05dc\u05de\u05ea \u05de\u05e1 \u05d5\u05dc\u05e4\u05e2\u05d9\u05dc\u05d5\u05d9\u05d5\u05ea \u05db\u05dc\u05db\u05dc\u05d9\u05d5\u05ea \u05d1\u05dc\u05ea\u05d9 \u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05d5\u05ea \\\"\u05dc\u05d4\u05d2\u05d9\u05d7\\\" \u05d1\u05d0\u05d5\u05e4\u05df \u05d8\u05d1\u05e2\u05d9 \u05de\u05ea\u05d5\u05da \u05d4\u05d0\u05d9\u05e0\u05d8\u05e8\u05d0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d1\u05d9\u05df \u05d4\u05e1\u05d5\u05db\u05e0\u05d9\u05dd, \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d5\u05d2\u05d3\u05e8\u05d5\u05ea \u05de\u05e8\u05d0\u05e9 \u05db\u05db\u05dc\u05dc\u05d9\u05dd \u05e7\u05e9\u05d9\u05d7\u05d9\u05dd. \u05de\u05d4 \u05d2\u05dd, \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05e8\u05d0\u05d5\u05ea \u05e9\u05d9\u05e0\u05d5\u05d9 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9 \u05d1\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d4\u05e8\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9\u05ea (DRL) \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e9\u05d9\u05e0\u05d5\u05d9 \u05de\u05e1\u05e4\u05d9\u05e7 \u05d0\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9 \u05e9\u05dc \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4 LLM \u05d1\u05e2\u05d6\u05e8\u05ea \u05e9\u05d9\u05e0\u05d5\u05d9\u05dd \u05db\u05de\u05d5 \u05ea\u05d9\u05d0\u05d5\u05e8 \u05d4\u05d0\u05d5\u05e4\u05d9 \u05e9\u05dc \u05d4\u05e1\u05d5\u05db\u05df. \",\n  \"\u05d2\u05dd \u05d0\u05dd \u05d0\u05ea\u05dd \u05dc\u05d0 \u05d7\u05d5\u05d1\u05d1\u05d9 \u05db\u05dc\u05db\u05dc\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd, \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d6\u05d0\u05ea \u05e9\u05dc \u05e9\u05d9\u05dc\u05d5\u05d1 \u05d1\u05d9\u05df LLM \u05dc DRL \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9LLM \u05de\u05e9\u05e8\u05ea \u05d0\u05ea \u05d4 DRL \u05d5\u05dc\u05d0 \u05d4\u05e4\u05d5\u05da (\u05db\u05de\u05d5 \u05e9\u05e7\u05d5\u05e8\u05d4 \u05d1\u05d0\u05d9\u05de\u05d5\u05df conversational LLMs \u05d0\u05d5 \u05e9\u05d0\u05ea\u05dd \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05d9\u05d6\u05d5 \u05ea\u05e9\u05d5\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d0\u05d4\u05d1\u05ea\u05dd \u05e9\u05dc chatGPT) \u05e4\u05d5\u05ea\u05d7\u05ea \u05d0\u05ea \u05d4\u05d3\u05dc\u05ea \u05dc\u05db\u05dc \u05de\u05d9\u05e0\u05d9 \u05e9\u05d9\u05de\u05d5\u05e9\u05d9\u05dd \u05d0\u05e4\u05dc\u05d9\u05e7\u05d8\u05d9\u05d1\u05d9\u05dd \u05e9\u05dc\u05d0 \u05d4\u05d9\u05d5 \u05db\u05dc \u05db\u05da \u05e0\u05d2\u05d9\u05e9\u05d9\u05dd \u05dc\u05e4\u05e0\u05d9 \u05d6\u05d4, \u05db\u05de\u05d5:\",\n  \"\u05d1\u05de\u05e7\u05d5\u05dd \u05e8\u05e7 \u05dc\u05d7\u05d6\u05d5\u05ea \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d1\u05d7\u05d9\u05e8\u05d5\u05ea, \u05d0\u05e4\u05e9\u05e8 \u05dc\u05d3\u05de\u05d5\u05ea \u05d0\u05d9\u05da \u05d3\u05e2\u05d5\u05ea \u05de\u05ea\u05e4\u05e9\u05d8\u05d5\u05ea, \u05d0\u05d9\u05da \u05e7\u05d1\u05d5\u05e6\u05d5\u05ea \u05d7\u05d1\u05e8\u05ea\u05d9\u05d5\u05ea \u05e0\u05d5\u05e6\u05e8\u05d5\u05ea \u05d0\u05d5 \u05de\u05ea\u05e4\u05e8\u05e7\u05d5\u05ea, \u05d0\u05d5 \u05d0\u05d9\u05da \u05de\u05ea\u05e4\u05ea\u05d7\u05ea \u05e7\u05d9\u05e6\u05d5\u05e0\u05d9\u05d5\u05ea \u2013 \u05dc\u05d0 \u05de\u05ea\u05d5\u05da \u05db\u05dc\u05dc\u05d9 \u05d1\u05e8\u05d5\u05e8 \u05d0\u05dc\u05d0 \u05de\u05d0\u05d9\u05e0\u05d8\u05e8\u05d0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea. \u05d0\u05e4\u05e9\u05e8 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05d9\u05da \u05e7\u05de\u05e4\u05d9\u05d9\u05df \u05de\u05e1\u05d5\u05d9\u05dd \u05d0\u05d5 \u05d7\u05d5\u05e7 \u05d7\u05d3\u05e9 \u05d9\u05e9\u05e4\u05d9\u05e2 \u05e2\u05dc \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d0\u05d6\u05e8\u05d7\u05d9\u05dd.\",\n  \"\u05d0\u05d9\u05da \u05e9\u05d9\u05e0\u05d5\u05d9 \u05d1\u05e0\u05ea\u05d9\u05d1 \u05ea\u05d7\u05d1\u05d5\u05e8\u05d4 \u05e6\u05d9\u05d1\u05d5\u05e8\u05d9\u05ea \u05d0\u05d5 \u05d1\u05e0\u05d9\u05d9\u05ea \u05e9\u05db\u05d5\u05e0\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05ea\u05e9\u05e4\u05d9\u05e2 \u05e2\u05dc \u05d3\u05e4\u05d5\u05e1\u05d9 \u05e0\u05e1\u05d9\u05e2\u05d4, \u05e4\u05e7\u05e7\u05d9\u05dd, \u05d0\u05d5 \u05d0\u05e4\u05d9\u05dc\u05d5 \u05e2\u05dc \u05e4\u05d9\u05ea\u05d5\u05d7 \u05e2\u05e1\u05e7\u05d9\u05dd \u05d1\u05d0\u05d6\u05d5\u05e8\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd, \u05d1\u05d2\u05dc\u05dc \u05d4\u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05d4\u05d3\u05d9\u05e0\u05de\u05d9\u05d5\u05ea \u05e9\u05dc \u05ea\u05d5\u05e9\u05d1\u05d9\u05dd \u05d5\u05e0\u05d4\u05d2\u05d9\u05dd.\",\n  \" \u05d0\u05d9\u05da \u05d7\u05d1\u05e8\u05d5\u05ea \u05de\u05d2\u05d9\u05d1\u05d5\u05ea \u05dc\u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05e9\u05dc \u05de\u05ea\u05d7\u05e8\u05d9\u05dd? \u05d4\u05d0\u05dd \u05d4\u05df \u05de\u05ea\u05db\u05e0\u05e1\u05d5\u05ea \u05dc\u05e7\u05e8\u05d0\u05ea \u05e7\u05e8\u05d8\u05dc \u05d0\u05d5 \u05e0\u05db\u05e0\u05e1\u05d5\u05ea \u05dc\u05de\u05dc\u05d7\u05de\u05ea \u05de\u05d7\u05d9\u05e8\u05d9\u05dd? \u05d0\u05e4\u05e9\u05e8 \u05dc\u05d3\u05de\u05d5\u05ea \u05d0\u05ea \u05d4\u05e9\u05d5\u05e7 \u05e2\u05dd \u05d7\u05d1\u05e8\u05d5\u05ea \\\"\u05d7\u05db\u05de\u05d5\u05ea\\\" \u05e9\u05de\u05e7\u05d1\u05dc\u05d5\u05ea \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05d0\u05e1\u05d8\u05e8\u05d8\u05d2\u05d9\u05d5\u05ea \u05d5\u05dc\u05e8\u05d0\u05d5\u05ea \u05de\u05d4\u05df \u05d4\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05d9\u05d5\u05ea \u05d4\u05e2\u05e1\u05e7\u05d9\u05d5\u05ea \u05d4\u05de\u05d2\u05d9\u05d7\u05d5\u05ea.\",\n  \"\u05d1\u05e7\u05d9\u05e6\u05d5\u05e8, \u05d6\u05d4 \u05dc\u05d0 \u05e8\u05e7 \u05e2\u05dc \u05d4\u05e2\u05dc\u05de\u05ea \u05de\u05e1. \u05d6\u05d5 \u05d3\u05e8\u05da \u05d7\u05d3\u05e9\u05d4 \u05d5\u05d9\u05e2\u05d9\u05dc\u05d4 \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05dc\u05db\u05dc \u05de\u05e2\u05e8\u05db\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05ea \u05e9\u05d1\u05d4 \u05d4\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d4\u05db\u05d5\u05dc\u05dc\u05ea \u05d4\u05d9\u05d0 \u05d9\u05d5\u05ea\u05e8 \u05de\u05e1\u05db\u05d5\u05dd \u05d7\u05dc\u05e7\u05d9\u05d4, \u05d5\u05de\u05d5\u05e9\u05e4\u05e2\u05ea \u05de\u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05d3\u05d9\u05e0\u05de\u05d9\u05d5\u05ea \u05d5\u05dc\u05de\u05d9\u05d3\u05d4 \u05e9\u05dc \u05d4\u05e4\u05e8\u05d8\u05d9\u05dd \u05d1\u05ea\u05d5\u05db\u05d4. \u05d6\u05d4 \u05e0\u05d5\u05ea\u05df \u05dc\u05e0\u05d5 \u05d9\u05db\u05d5\u05dc\u05ea \\\"\u05dc\u05e9\u05d7\u05e7\\\" \u05e2\u05dd \u05d4\u05de\u05e6\u05d9\u05d0\u05d5\u05ea, \u05dc\u05d1\u05d3\u05d5\u05e7 \u05ea\u05e8\u05d7\u05d9\u05e9\u05d9\u05dd \u05d5\u05dc\u05dc\u05de\u05d5\u05d3 \u05de\u05d4\u05dd, \u05d1\u05dc\u05d9 \u05d4\u05e6\u05d5\u05e8\u05da \u05dc\u05ea\u05db\u05e0\u05ea \u05de\u05e8\u05d0\u05e9 \u05db\u05dc \u05e4\u05e8\u05d8.\",\n  \"\u05dc\u05d0 \u05de\u05d0\u05de\u05e8 \u05e7\u05dc\u05d0\u05e1\u05d9 \u05de\u05de\u05d4 \u05e9\u05e2\u05d5\u05dc\u05d4 \u05e4\u05d4 \u05d1\u05e1\u05e7\u05d9\u05e8\u05d4 \u05d1\u05d3\u05e8\u05da \u05db\u05dc\u05dc, \u05d0\u05d1\u05dc \u05d9\u05db\u05d5\u05dc \u05dc\u05e4\u05ea\u05d5\u05d7 \u05d0\u05ea \u05d4\u05e8\u05d0\u05e9:\",\n  \"https://arxiv.org/abs/2501.18177\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst existingCount = paragraphs.items.length; // 13 in before.docx\n\n// 1) Overwrite the text of each existing paragraph (0-12) with the new\n//    review's corresponding paragraph text. Using insertText(..., \"Replace\")\n//    on the paragraph's own range rewrites its run(s) in place - for\n//    paragraph 0 this also drops the old <w:br/> + second run, since the\n//    whole paragraph range (including the line break) is replaced.\nfor (let i = 0; i < existingCount; i++) {\n  paragraphs.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Append the two new trailing paragraphs that didn't exist before:\n//    - a closing remark paragraph\n//    - a paragraph holding a line break followed by the new arXiv link\nbody.insertParagraph(newTexts[existingCount], Word.InsertLocation.end);\nawait context.sync();\n\n// Insert as a single run containing both the break and the url text by\n// prefixing the text with the vertical-tab line-break character (\\u000B),\n// which Word serializes as <w:br/> inside the same run.\nbody.insertParagraph(\"\\u000b\" + newTexts[existingCount + 1], Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# PowerShell / Word COM-interop edit script.\n# Applies the Review_480 edit: updates the date, swaps the paper title and\n# review body paragraphs from the \"Procedural Knowledge\" review to the\n# \"Tax Evasion\" review, and appends a closing remark + new arXiv link\n# (preceded by a line break), replacing the old link paragraph.\n\n$d = $word.ActiveDocument\n\n$newTexts = @(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7: 04.07.25\",\n  \"Investigating Tax Evasion Emergence Using Dual Large Language Model and Deep Reinforcement Learning Powered Agent-based Simulation\",\n  \"\u05ea\u05e4\u05e0\u05d9\u05ea \u05de\u05e4\u05ea\u05d9\u05e2\u05d4 \u05de\u05ea\u05e8\u05d7\u05e9\u05ea \u05d1\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1-LLMs \u05d1\u05ea\u05d7\u05d5\u05de\u05d9\u05dd `\"\u05e8\u05db\u05d9\u05dd`\" \u05d9\u05d5\u05ea\u05e8 \u05db\u05de\u05d5 \u05e4\u05e1\u05d9\u05db\u05d5\u05dc\u05d5\u05d2\u05d9\u05d4, \u05e1\u05d5\u05e6\u05d9\u05d5\u05dc\u05d5\u05d2\u05d9\u05d4 \u05d5\u05d0\u05e4\u05d9\u05dc\u05d5 \u05db\u05dc\u05db\u05dc\u05d4. \u05d0\u05d5\u05de\u05e0\u05dd LLMs \u05dc\u05d0 `\"\u05d7\u05d5\u05e9\u05d1\u05d9\u05dd`\" \u05db\u05de\u05d5 \u05d1\u05e0\u05d9 \u05d0\u05d3\u05dd \u05d1\u05e8\u05de\u05ea \u05d4\u05d0\u05d9\u05e0\u05d3\u05d9\u05d1\u05d9\u05d3\u05d5\u05d0\u05dc, \u05d0\u05d1\u05dc \u05de\u05e1\u05ea\u05d1\u05e8 \u05e9\u05d4\u05dd \u05db\u05d1\u05e8 \u05de\u05d7\u05db\u05d9\u05dd \u05d0\u05ea \u05d0\u05d9\u05da \u05e9\u05d0\u05e0\u05d7\u05e0\u05d5 \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05db\u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d4. \",\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05e0\u05d5\u05db\u05d7\u05d9 \u05de\u05e6\u05d9\u05d2 \u05d2\u05d9\u05e9\u05d4 \u05d7\u05d3\u05e9\u05e0\u05d9\u05ea \u05dc\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1 LLMs \u05dc\u05d7\u05e7\u05e8 \u05d4\u05ea\u05d7\u05de\u05e7\u05d5\u05ea \u05de\u05de\u05e1 \u05d1\u05e2\u05d6\u05e8\u05ea \u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4. \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05d4\u05e0\u05d9\u05d7 \u05de\u05e8\u05d0\u05e9 \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05e9\u05dc \u05d4\u05e2\u05dc\u05de\u05ea \u05de\u05e1, \u05db\u05de\u05d5 \u05e9\u05e2\u05e9\u05d5 \u05db\u05dc \u05d4\u05de\u05d7\u05e7\u05e8\u05d9\u05dd \u05e2\u05d3 \u05d4\u05d9\u05d5\u05dd, \u05d4\u05de\u05d7\u05e7\u05e8 \u05de\u05ea\u05de\u05e7\u05d3 \u05d1\u05d4\u05d5\u05e4\u05e2\u05d4 \u05d5\u05d1\u05d3\u05d9\u05e0\u05de\u05d9\u05e7\u05d4 \u05e9\u05dc \u05ea\u05d5\u05e4\u05e2\u05d4 \u05d6\u05d5 \u05d1\u05e7\u05e8\u05d1 \u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d9\u05d4. \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e1\u05d5\u05db\u05e0\u05d9\u05dd (agent-based simulation) \u05d4\u05de\u05e9\u05dc\u05d1\u05ea LLMs \u05d5\u05dc\u05de\u05d9\u05d3\u05d4 \u05e2\u05de\u05d5\u05e7\u05d4 \u05e2\u05dd \u05d7\u05d9\u05d6\u05d5\u05e7\u05d9\u05dd (deep reinforcement learning), \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05d1\u05d5\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05d4\u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05d9\u05d5\u05ea \u05db\u05dc\u05db\u05dc\u05d9\u05d5\u05ea \u05d1\u05dc\u05ea\u05d9 \u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05d5\u05ea (\u05de\u05d4 \u05e9\u05d4\u05e8\u05d1 \u05de\u05db\u05d9\u05e8\u05d9\u05dd \u05d1\u05ea\u05d5\u05e8 `\"\u05db\u05dc\u05db\u05dc\u05d4 \u05e9\u05d7\u05d5\u05e8\u05d4`\") \u05dc\u05d4\u05d2\u05d9\u05d7 \u05d1\u05d0\u05d5\u05e4\u05df \u05e1\u05e4\u05d5\u05e0\u05d8\u05e0\u05d9, \u05d5\u05dc\u05d0 \u05db\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05de\u05d5\u05d2\u05d3\u05e8\u05ea \u05de\u05e8\u05d0\u05e9. \u05e2\u05dc \u05d0\u05e3 \u05d4\u05ea\u05e8\u05d5\u05de\u05d4 \u05e9\u05dc\u05d5 \u05dc\u05db\u05dc\u05db\u05dc\u05d4, \u05de\u05d4 \u05e9\u05dc\u05db\u05e0\u05e8\u05d0\u05d4 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d0\u05d5\u05ea\u05e0\u05d5 \u05d6\u05d4 \u05d4\u05e9\u05d9\u05de\u05d5\u05e9 \u05d4\u05d9\u05d9\u05d7\u05d5\u05d3\u05d9 \u05d1 LLMs \u05d5\u05d1 DRL \u05db\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc \u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05e9\u05d9\u05d5\u05d3\u05e2 \u05dc\u05e7\u05d7\u05ea \u05de\u05d9\u05d3\u05e2 \u05dc\u05d0 \u05e4\u05d5\u05e8\u05de\u05dc\u05d9 (\u05db\u05de\u05d5 \u05ea\u05d9\u05d0\u05d5\u05e8 \u05d0\u05d9\u05e9\u05d9\u05d5\u05ea) \u05d5\u05dc\u05d4\u05db\u05e0\u05d9\u05e1 \u05d0\u05d5\u05ea\u05d5 \u05dc\u05e1\u05d8 \u05d4\u05e9\u05d9\u05e7\u05d5\u05dc\u05d9\u05dd \u05d4\u05e4\u05d5\u05e8\u05de\u05dc\u05d9 \u05e9\u05dc \u05e1\u05d5\u05db\u05df - \u05dc\u05de\u05e9\u05dc \u05db\u05de\u05d4 \u05d0\u05d5\u05e4\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc \u05e4\u05e8\u05e1\u05d5\u05e0\u05dc\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d1\u05de\u05e2\u05e8\u05db\u05d5\u05ea \u05d0\u05e4\u05e9\u05e8 \u05dc\u05e2\u05e9\u05d5\u05ea \u05e2\u05dc \u05d2\u05d1\u05d9 \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05d6\u05d4.\",\n  \"\u05d4\u05e9\u05d9\u05d8\u05d4 \u05d1\u05de\u05d0\u05de\u05e8 \u05de\u05ea\u05e2\u05de\u05e7\u05ea \u05d1\u05d9\u05e6\u05d9\u05e8\u05ea \u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05d4\u05de\u05d3\u05de\u05d4 \u05db\u05dc\u05db\u05dc\u05d4 \u05e1\u05d2\u05d5\u05e8\u05d4 (\u05d9\u05e9 \u05de\u05e1\u05d7\u05e8 \u05d1\u05d9\u05df \u05d0\u05e0\u05e9\u05d9\u05dd \u05d1\u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d4 \u05e9\u05de\u05e0\u05e1\u05d4 \u05dc\u05d4\u05d8\u05d9\u05d9\u05d1 \u05e2\u05dd \u05de\u05e6\u05d1\u05dd), \u05e9\u05d1\u05d4 \u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05e4\u05d5\u05e2\u05dc\u05d9\u05dd \u05d5\u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea. \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05d4\u05e9\u05ea\u05de\u05e9\u05d5 \u05d1\u05de\u05d1\u05e0\u05d4 \u05db\u05dc\u05db\u05dc\u05d4 \u05d3\u05d5\u05de\u05d4 \u05dc\u05d6\u05d0\u05ea \u05e9\u05dc \u05d0\u05e8\u05e6\u05d5\u05ea \u05d4\u05d1\u05e8\u05d9\u05ea \u05e2\u05dd \u05d3\u05d9\u05d5\u05d5\u05d7 \u05de\u05e1 \u05e2\u05e6\u05de\u05d9 \u05db\u05d3\u05d9 \u05dc\u05d0\u05e4\u05e9\u05e8 \u05dc\u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05d4\u05d6\u05d3\u05de\u05e0\u05d5\u05ea \u05dc\u05d4\u05e2\u05dc\u05d9\u05dd \u05de\u05e1 \u05db\u05d4\u05d7\u05dc\u05d8\u05d4. \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05d9\u05e9 \u05d4\u05de\u05d5\u05df \u05e1\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05de\u05e1 \u05d5\u05d6\u05d4 \u05de\u05e1\u05ea\u05d1\u05da \u05de\u05d4\u05e8 (\u05ea\u05e9\u05d0\u05dc\u05d5 \u05d0\u05ea \u05e8\u05d5\u05d0\u05d4 \u05d4\u05d7\u05e9\u05d1\u05d5\u05df \u05e9\u05dc\u05db\u05dd) \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05d4\u05ea\u05e8\u05db\u05d6\u05d5 \u05d1\u05e9\u05e0\u05d9 \u05e1\u05d5\u05d2\u05d9 \u05de\u05e1 - \u05de\u05e1 \u05d4\u05db\u05e0\u05e1\u05d4 \u05d5\u05de\u05e1 \u05e2\u05e8\u05da \u05de\u05d5\u05e1\u05e3 (\u05de\u05e2`\"\u05de). \",\n  \"\u05d4\u05dd \u05d2\u05dd \u05e1\u05d9\u05de\u05dc\u05e6\u05d5 \u05e8\u05e9\u05d5\u05d9\u05d5\u05ea \u05d4\u05d7\u05d5\u05e7 \u05d5\u05ea\u05d5\u05e2\u05dc\u05ea \u05e9\u05d4\u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05de\u05d4\u05de\u05d3\u05d9\u05e0\u05d4 \u05de\u05d4\u05de\u05e1 \u05e9\u05d4\u05dd \u05de\u05e9\u05dc\u05de\u05d9\u05dd \u05dc\u05d4. \u05d4\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05e2\u05e6\u05de\u05d4 \u05de\u05d4\u05d5\u05d5\u05d4 \u05ea\u05e9\u05ea\u05d9\u05ea \u05dc\u05d7\u05dc\u05e7 \u05d4\u05de\u05e8\u05db\u05d6\u05d9 \u05e9\u05dc \u05d4\u05e2\u05d1\u05d5\u05d3\u05d4 - \u05de\u05d5\u05d3\u05dc \u05e7\u05d1\u05dc\u05ea \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05e9\u05dc \u05d4\u05e1\u05d5\u05db\u05e0\u05d9\u05dd \u05d1\u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d4. \u05db\u05d3\u05d9 \u05dc\u05d9\u05d9\u05e6\u05e8 \u05d0\u05d5\u05db\u05dc\u05d5\u05e1\u05d9\u05d4 \u05d4\u05d8\u05e8\u05d5\u05d2\u05e0\u05d9\u05ea, \u05de\u05d5\u05d7 \u05e9\u05dc \u05e1\u05d5\u05db\u05df \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc LLM \u05d5 DRL \u05db\u05d0\u05e9\u05e8 LLMs \u05e9\u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05ea\u05d9\u05d0\u05d5\u05e8 \u05e9\u05dc \u05d4\u05d0\u05d5\u05e4\u05d9 \u05e9\u05dc \u05d4\u05e1\u05d5\u05db\u05df (\u05de\u05d1\u05d5\u05e1\u05e1 \u05e2\u05dc \u05d4\u05ea\u05d5\u05db\u05df \u05e9\u05d4\u05d5\u05d0 \u05de\u05e4\u05e8\u05e1\u05dd \u05d1\u05d8\u05d5\u05d5\u05d9\u05d8\u05e8 \u05dc\u05de\u05e9\u05dc), \u05d4\u05d9\u05e1\u05d8\u05d5\u05e8\u05d9\u05d9\u05ea \u05d4\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05db\u05d8\u05e7\u05e1\u05d8, \u05d5\u05d0\u05ea \u05db\u05dc \u05d4\u05de\u05d9\u05d3\u05e2 \u05e2\u05dc \u05d4\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05db\u05d8\u05e7\u05e1\u05d8 \u05d2\u05dd \u05db\u05df. \",\n  \"\u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4-context prompt \u05d4\u05d6\u05d4 \u05d4 LLM \u05e0\u05e9\u05d0\u05dc `\"\u05db\u05de\u05d4 \u05de\u05e1 \u05d0\u05e0\u05d9 \u05e6\u05e8\u05d9\u05da \u05dc\u05e9\u05dc\u05dd?`\". \u05d4\u05de\u05e1\u05e4\u05e8 \u05e9\u05d4 LLM \u05de\u05d7\u05d6\u05d9\u05e8, \u05de\u05d5\u05e2\u05d1\u05e8 \u05db\u05e7\u05dc\u05d8 \u05dc\u05de\u05d5\u05d3\u05dc  DRL \u05e9\u05de\u05e7\u05d1\u05dc \u05d2\u05dd \u05d0\u05ea \u05e9\u05d0\u05e8 \u05d4\u05d3\u05d0\u05d8\u05d4 \u05e9\u05e7\u05d9\u05d1\u05dc \u05d4 LLM \u05d0\u05d1\u05dc \u05d2\u05dd \u05db\u05de\u05d4 \u05d4\u05e1\u05d5\u05db\u05df `\"\u05d4\u05e8\u05e4\u05ea\u05e7\u05e0\u05d9`\" \u05db\u05e4\u05e8\u05de\u05d8\u05e8 \u05e9\u05d4-DRL \u05de\u05e9\u05ea\u05de\u05e9 \u05db\u05d3\u05d9 \u05dc\u05e2\u05e9\u05d5\u05ea \u05d0\u05e7\u05e1\u05e4\u05dc\u05d5\u05e8\u05e6\u05d9\u05d4. \u05db\u05dc\u05d5\u05de\u05e8, \u05d4 LLM \u05de\u05d7\u05d6\u05d9\u05e8 \u05d4\u05d7\u05dc\u05d8\u05d4 \u05e8\u05d0\u05e9\u05d5\u05e0\u05d9\u05ea \u05e9\u05d0\u05d5\u05ea\u05d4, \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05dc\u05e7\u05d8 \u05e9\u05dc \u05d4 LLM \u05d5\u05e2\u05d5\u05d3 \u05de\u05e9\u05ea\u05e0\u05d4 exploration (\u05d0\u05d4\u05d1\u05ea \u05e1\u05d9\u05db\u05d5\u05df \u05d1\u05de\u05d5\u05d1\u05df \u05d4\u05db\u05dc\u05db\u05dc\u05d9) \u05de\u05e7\u05d1\u05dc \u05d2\u05dd DRL \u05e9\u05de\u05e7\u05d1\u05dc \u05d4\u05d7\u05dc\u05d8\u05d4 \u05d1\u05e2\u05e6\u05de\u05d5 \u05e9\u05d4\u05d9\u05e0\u05d4 \u05d2\u05dd \u05d4\u05e1\u05d5\u05e4\u05d9\u05ea. \",\n  \"\u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05dc\u05de\u05d9\u05d3\u05d4 \u05d4\u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05e2\u05dc\u05de\u05ea \u05de\u05e1 \u05d5\u05dc\u05e4\u05e2\u05d9\u05dc\u05d5\u05d9\u05d5\u05ea \u05db\u05dc\u05db\u05dc\u05d9\u05d5\u05ea \u05d1\u05dc\u05ea\u05d9 \u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05d5\u05ea `\"\u05dc\u05d4\u05d2\u05d9\u05d7`\" \u05d1\u05d0\u05d5\u05e4\u05df \u05d8\u05d1\u05e2\u05d9 \u05de\u05ea\u05d5\u05da \u05d4\u05d0\u05d9\u05e0\u05d8\u05e8\u05d0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d1\u05d9\u05df \u05d4\u05e1\u05d5\u05db\u05e0\u05d9\u05dd, \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d5\u05d2\u05d3\u05e8\u05d5\u05ea \u05de\u05e8\u05d0\u05e9 \u05db\u05db\u05dc\u05dc\u05d9\u05dd \u05e7\u05e9\u05d9\u05d7\u05d9\u05dd. \u05de\u05d4 \u05d2\u05dd, \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05e8\u05d0\u05d5\u05ea \u05e9\u05d9\u05e0\u05d5\u05d9 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9 \u05d1\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d4\u05e8\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9\u05ea (DRL) \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e9\u05d9\u05e0\u05d5\u05d9 \u05de\u05e1\u05e4\u05d9\u05e7 \u05d0\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9 \u05e9\u05dc \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4 LLM \u05d1\u05e2\u05d6\u05e8\u05ea \u05e9\u05d9\u05e0\u05d5\u05d9\u05dd \u05db\u05de\u05d5 \u05ea\u05d9\u05d0\u05d5\u05e8 \u05d4\u05d0\u05d5\u05e4\u05d9 \u05e9\u05dc \u05d4\u05e1\u05d5\u05db\u05df. \",\n  \"\u05d2\u05dd \u05d0\u05dd \u05d0\u05ea\u05dd \u05dc\u05d0 \u05d7\u05d5\u05d1\u05d1\u05d9 \u05db\u05dc\u05db\u05dc\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd, \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d6\u05d0\u05ea \u05e9\u05dc \u05e9\u05d9\u05dc\u05d5\u05d1 \u05d1\u05d9\u05df LLM \u05dc DRL \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9LLM \u05de\u05e9\u05e8\u05ea \u05d0\u05ea \u05d4 DRL \u05d5\u05dc\u05d0 \u05d4\u05e4\u05d5\u05da (\u05db\u05de\u05d5 \u05e9\u05e7\u05d5\u05e8\u05d4 \u05d1\u05d0\u05d9\u05de\u05d5\u05df conversational LLMs \u05d0\u05d5 \u05e9\u05d0\u05ea\u05dd \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05d9\u05d6\u05d5 \u05ea\u05e9\u05d5\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d0\u05d4\u05d1\u05ea\u05dd \u05e9\u05dc chatGPT) \u05e4\u05d5\u05ea\u05d7\u05ea \u05d0\u05ea \u05d4\u05d3\u05dc\u05ea \u05dc\u05db\u05dc \u05de\u05d9\u05e0\u05d9 \u05e9\u05d9\u05de\u05d5\u05e9\u05d9\u05dd \u05d0\u05e4\u05dc\u05d9\u05e7\u05d8\u05d9\u05d1\u05d9\u05dd \u05e9\u05dc\u05d0 \u05d4\u05d9\u05d5 \u05db\u05dc \u05db\u05da \u05e0\u05d2\u05d9\u05e9\u05d9\u05dd \u05dc\u05e4\u05e0\u05d9 \u05d6\u05d4, \u05db\u05de\u05d5:\",\n  \"\u05d1\u05de\u05e7\u05d5\u05dd \u05e8\u05e7 \u05dc\u05d7\u05d6\u05d5\u05ea \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d1\u05d7\u05d9\u05e8\u05d5\u05ea, \u05d0\u05e4\u05e9\u05e8 \u05dc\u05d3\u05de\u05d5\u05ea \u05d0\u05d9\u05da \u05d3\u05e2\u05d5\u05ea \u05de\u05ea\u05e4\u05e9\u05d8\u05d5\u05ea, \u05d0\u05d9\u05da \u05e7\u05d1\u05d5\u05e6\u05d5\u05ea \u05d7\u05d1\u05e8\u05ea\u05d9\u05d5\u05ea \u05e0\u05d5\u05e6\u05e8\u05d5\u05ea \u05d0\u05d5 \u05de\u05ea\u05e4\u05e8\u05e7\u05d5\u05ea, \u05d0\u05d5 \u05d0\u05d9\u05da \u05de\u05ea\u05e4\u05ea\u05d7\u05ea \u05e7\u05d9\u05e6\u05d5\u05e0\u05d9\u05d5\u05ea \u2013 \u05dc\u05d0 \u05de\u05ea\u05d5\u05da \u05db\u05dc\u05dc\u05d9 \u05d1\u05e8\u05d5\u05e8 \u05d0\u05dc\u05d0 \u05de\u05d0\u05d9\u05e0\u05d8\u05e8\u05d0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea. \u05d0\u05e4\u05e9\u05e8 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05d9\u05da \u05e7\u05de\u05e4\u05d9\u05d9\u05df \u05de\u05e1\u05d5\u05d9\u05dd \u05d0\u05d5 \u05d7\u05d5\u05e7 \u05d7\u05d3\u05e9 \u05d9\u05e9\u05e4\u05d9\u05e2 \u05e2\u05dc \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d0\u05d6\u05e8\u05d7\u05d9\u05dd.\",\n  \"\u05d0\u05d9\u05da \u05e9\u05d9\u05e0\u05d5\u05d9 \u05d1\u05e0\u05ea\u05d9\u05d1 \u05ea\u05d7\u05d1\u05d5\u05e8\u05d4 \u05e6\u05d9\u05d1\u05d5\u05e8\u05d9\u05ea \u05d0\u05d5 \u05d1\u05e0\u05d9\u05d9\u05ea \u05e9\u05db\u05d5\u05e0\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05ea\u05e9\u05e4\u05d9\u05e2 \u05e2\u05dc \u05d3\u05e4\u05d5\u05e1\u05d9 \u05e0\u05e1\u05d9\u05e2\u05d4, \u05e4\u05e7\u05e7\u05d9\u05dd, \u05d0\u05d5 \u05d0\u05e4\u05d9\u05dc\u05d5 \u05e2\u05dc \u05e4\u05d9\u05ea\u05d5\u05d7 \u05e2\u05e1\u05e7\u05d9\u05dd \u05d1\u05d0\u05d6\u05d5\u05e8\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd, \u05d1\u05d2\u05dc\u05dc \u05d4\u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05d4\u05d3\u05d9\u05e0\u05de\u05d9\u05d5\u05ea \u05e9\u05dc \u05ea\u05d5\u05e9\u05d1\u05d9\u05dd \u05d5\u05e0\u05d4\u05d2\u05d9\u05dd.\",\n  \" \u05d0\u05d9\u05da \u05d7\u05d1\u05e8\u05d5\u05ea \u05de\u05d2\u05d9\u05d1\u05d5\u05ea \u05dc\u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05e9\u05dc \u05de\u05ea\u05d7\u05e8\u05d9\u05dd? \u05d4\u05d0\u05dd \u05d4\u05df \u05de\u05ea\u05db\u05e0\u05e1\u05d5\u05ea \u05dc\u05e7\u05e8\u05d0\u05ea \u05e7\u05e8\u05d8\u05dc \u05d0\u05d5 \u05e0\u05db\u05e0\u05e1\u05d5\u05ea \u05dc\u05de\u05dc\u05d7\u05de\u05ea \u05de\u05d7\u05d9\u05e8\u05d9\u05dd? \u05d0\u05e4\u05e9\u05e8 \u05dc\u05d3\u05de\u05d5\u05ea \u05d0\u05ea \u05d4\u05e9\u05d5\u05e7 \u05e2\u05dd \u05d7\u05d1\u05e8\u05d5\u05ea `\"\u05d7\u05db\u05de\u05d5\u05ea`\" \u05e9\u05de\u05e7\u05d1\u05dc\u05d5\u05ea \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05d0\u05e1\u05d8\u05e8\u05d8\u05d2\u05d9\u05d5\u05ea \u05d5\u05dc\u05e8\u05d0\u05d5\u05ea \u05de\u05d4\u05df \u05d4\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05d9\u05d5\u05ea \u05d4\u05e2\u05e1\u05e7\u05d9\u05d5\u05ea \u05d4\u05de\u05d2\u05d9\u05d7\u05d5\u05ea.\",\n  \"\u05d1\u05e7\u05d9\u05e6\u05d5\u05e8, \u05d6\u05d4 \u05dc\u05d0 \u05e8\u05e7 \u05e2\u05dc \u05d4\u05e2\u05dc\u05de\u05ea \u05de\u05e1. \u05d6\u05d5 \u05d3\u05e8\u05da \u05d7\u05d3\u05e9\u05d4 \u05d5\u05d9\u05e2\u05d9\u05dc\u05d4 \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05dc\u05db\u05dc \u05de\u05e2\u05e8\u05db\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05ea \u05e9\u05d1\u05d4 \u05d4\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d4\u05db\u05d5\u05dc\u05dc\u05ea \u05d4\u05d9\u05d0 \u05d9\u05d5\u05ea\u05e8 \u05de\u05e1\u05db\u05d5\u05dd \u05d7\u05dc\u05e7\u05d9\u05d4, \u05d5\u05de\u05d5\u05e9\u05e4\u05e2\u05ea \u05de\u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05d3\u05d9\u05e0\u05de\u05d9\u05d5\u05ea \u05d5\u05dc\u05de\u05d9\u05d3\u05d4 \u05e9\u05dc \u05d4\u05e4\u05e8\u05d8\u05d9\u05dd \u05d1\u05ea\u05d5\u05db\u05d4. \u05d6\u05d4 \u05e0\u05d5\u05ea\u05df \u05dc\u05e0\u05d5 \u05d9\u05db\u05d5\u05dc\u05ea `\"\u05dc\u05e9\u05d7\u05e7`\" \u05e2\u05dd \u05d4\u05de\u05e6\u05d9\u05d0\u05d5\u05ea, \u05dc\u05d1\u05d3\u05d5\u05e7 \u05ea\u05e8\u05d7\u05d9\u05e9\u05d9\u05dd \u05d5\u05dc\u05dc\u05de\u05d5\u05d3 \u05de\u05d4\u05dd, \u05d1\u05dc\u05d9 \u05d4\u05e6\u05d5\u05e8\u05da \u05dc\u05ea\u05db\u05e0\u05ea \u05de\u05e8\u05d0\u05e9 \u05db\u05dc \u05e4\u05e8\u05d8.\",\n  \"\u05dc\u05d0 \u05de\u05d0\u05de\u05e8 \u05e7\u05dc\u05d0\u05e1\u05d9 \u05de\u05de\u05d4 \u05e9\u05e2\u05d5\u05dc\u05d4 \u05e4\u05d4 \u05d1\u05e1\u05e7\u05d9\u05e8\u05d4 \u05d1\u05d3\u05e8\u05da \u05db\u05dc\u05dc, \u05d0\u05d1\u05dc \u05d9\u05db\u05d5\u05dc \u05dc\u05e4\u05ea\u05d5\u05d7 \u05d0\u05ea \u05d4\u05e8\u05d0\u05e9:\",\n  \"https://arxiv.org/abs/2501.18177\"\n)\n\n$existingCount = $d.Paragraphs.Count  # 13 in before.docx\n\n# 1) Overwrite the text of each existing paragraph (1..13) with the new\n#    review's corresponding paragraph text. Assigning Range.Text rewrites\n#    the run(s) in place while keeping the paragraph mark - for paragraph 1\n#    this also drops the old line break + second run, since the whole\n#    paragraph range (including the line break) is replaced.\nfor ($i = 1; $i -le $existingCount; $i++) {\n  $d.Paragraphs.Item($i).Range.Text = $newTexts[$i - 1]\n}\n\n# 2) Append the two new trailing paragraphs that didn't exist before:\n#    - a closing remark paragraph\n#    - a paragraph holding a line break followed by the new arXiv link\n$lastRange = $d.Paragraphs.Item($existingCount).Range\n$lastRange.InsertParagraphAfter()\n$d.Paragraphs.Item($existingCount + 1).Range.Text = $newTexts[$existingCount]\n\n$secondLastRange = $d.Paragraphs.Item($existingCount + 1).Range\n$secondLastRange.InsertParagraphAfter()\n$finalRange = $d.Paragraphs.Item($existingCount + 2).Range\n$finalRange.InsertAfter([char]11)\n$finalRange.Collapse(0)\n$finalRange.InsertAfter($newTexts[$existingCount + 1])\n\nWrite-Output $d.Paragraphs.Count\n"}
